$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new export run (2024-09-24 08:58:15)
$ws.Name = "IClientBalance-20240924-085815-"

# Determine the extent of the data (header in row 1, data starting row 2)
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()

# Bump the "Dt. Referencia" column (G) by one day for every data row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cur = $cell.Value2()
    if ($cur -ne $null) {
        $cell.Value = $cur + 1
    }
}

# Correct a few data-entry typos in "Saldo Previsto" (E) / "Vl. Total" (H)
$ws.Cells.Item(112, 5).Value = 0.14
$ws.Cells.Item(112, 8).Value = 0.14

$ws.Cells.Item(113, 5).Value = 0.81
$ws.Cells.Item(113, 8).Value = 0.81

$ws.Cells.Item(118, 5).Value = 8916.85
$ws.Cells.Item(118, 8).Value = 8916.85
